$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected Price (D) / Volume(1h) (E) cells to remain text so
# Excel does not reinterpret the numeric-looking / percent-looking strings
# as Number/Percentage cells (source data is stored as plain text).
$ws.Range("D2:E47").NumberFormat = "@"

$ws.Range("D2").Value = "285.25"
$ws.Range("E2").Value = "2.11%"
$ws.Range("D3").Value = "28.72"
$ws.Range("E3").Value = "4.79%"
$ws.Range("D4").Value = "5.094"
$ws.Range("E4").Value = "5.37%"
$ws.Range("D5").Value = "0.06682"
$ws.Range("E5").Value = "4.71%"
$ws.Range("D6").Value = "7.333"
$ws.Range("E6").Value = "4.30%"
$ws.Range("D7").Value = "3.381"
$ws.Range("E7").Value = "2.46%"
$ws.Range("D8").Value = "1.348"
$ws.Range("E8").Value = "4.34%"
$ws.Range("D9").Value = "0.9343"
$ws.Range("E9").Value = "4.70%"
$ws.Range("D10").Value = "0.1571"
$ws.Range("E10").Value = "3.16%"
$ws.Range("D11").Value = "0.06594"
$ws.Range("E11").Value = "16.86%"
$ws.Range("D12").Value = "0.07683"
$ws.Range("E12").Value = "2.51%"
$ws.Range("D13").Value = "0.02905"
$ws.Range("E13").Value = "-0.41%"
$ws.Range("E14").Value = "-0.14%"
$ws.Range("D15").Value = "0.001588"
$ws.Range("E15").Value = "0.96%"
$ws.Range("D16").Value = "0.04471"
$ws.Range("E16").Value = "1.89%"
$ws.Range("D17").Value = "0.0006436"
$ws.Range("E17").Value = "0.79%"
$ws.Range("D18").Value = "0.006529"
$ws.Range("E18").Value = "6.84%"
$ws.Range("D19").Value = "3.485"
$ws.Range("E19").Value = "0.42%"
$ws.Range("D20").Value = "2.238"
$ws.Range("E20").Value = "-2.47%"
$ws.Range("E22").Value = "-3.32%"
$ws.Range("D23").Value = "4.048"
$ws.Range("E23").Value = "4.00%"
$ws.Range("E24").Value = "1.14%"
$ws.Range("D25").Value = "0.001176"
$ws.Range("E25").Value = "0.12%"
$ws.Range("D26").Value = "0.004475"
$ws.Range("E26").Value = "4.53%"
$ws.Range("D27").Value = "0.0001245"
$ws.Range("E27").Value = "5.82%"
$ws.Range("E28").Value = "-2.46%"
$ws.Range("D40").Value = "0.04199"
$ws.Range("E40").Value = "3.56%"
$ws.Range("D41").Value = "0.006731"
$ws.Range("E41").Value = "-0.13%"
$ws.Range("E42").Value = "-11.37%"
$ws.Range("D43").Value = "0.002012"
$ws.Range("E43").Value = "-1.56%"
$ws.Range("D44").Value = "0.01218"
$ws.Range("E44").Value = "8.97%"
$ws.Range("D45").Value = "0.00005662"
$ws.Range("E45").Value = "1.96%"
$ws.Range("E47").Value = "-29.55%"

# Restore the default "Normal" style on the touched range so no stray
# number-format style survives the text-coercion trick above.
$ws.Range("D2:E47").Style = "Normal"

